$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$v = $ws.Range("A1").Value2
Write-Host "A1: $v"
$v2 = $ws.Cells.Item(1,1).Text
Write-Host "A1 text: $v2"
